$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "point sum" value for Labs section
$ws.Range("D4").Value = 340.35

# Pre-Lab section (rows 12-17) - mark as complete
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 1

# Projects section (rows 28-33)
$ws.Range("D28").Value = 20
$ws.Range("D29").Value = 29
$ws.Range("D30").Value = 38
$ws.Range("D31").Value = 50
$ws.Range("D32").Value = 50
$ws.Range("D33").Value = 54

# Exams section (rows 42-43)
$ws.Range("D42").Value = 130
$ws.Range("D43").Value = 110

# Labs section (rows 57-58)
$ws.Range("D57").Value = 2.5
$ws.Range("D58").Value = 1.43

# Update selection to match new active cell
$ws.Range("I4").Select()
